$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.665.99"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "1.869.09"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.87"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +3.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4685"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +4.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3940"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +2.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.56"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08027"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +1.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.022"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("E12").Value = "  +2.16%  "

$ws.Range("D13").Value = "1.871.57"
$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.928"
$ws.Range("D14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.130"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001045"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.58"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06636"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.20"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +1.87%  "

$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "27.680.11"
$ws.Range("E22").Value = "  +1.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.484"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("E24").Value = "  +2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("D26").Value = "2.095.24"
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.34"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +4.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.18"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +2.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.089"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.547"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +2.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.60"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9635"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09477"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.447"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.594"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.308"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("E37").Value = "  +1.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06066"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +1.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.225"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +1.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.104"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.75%  "

$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("E42").Value = "  +1.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1892"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.24"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.254"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5677"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.16"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.384"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.934"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06846"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.32"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +5.87%  "
